# Update "Forecast Comparison" sheet with a new Week_Start_Date column,
# corrected week labels/types, and corrected MyForecast value; then
# update the dependent totals on the "Summary" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# 1. Insert a new column B ("Week_Start_Date") before the existing ASIN column.
$ws.Columns.Item(2).Insert()

# 2. Header for the new column.
$ws.Cells.Item(1, 2).Value2 = "Week_Start_Date"

# 3. Week_Start_Date values (weekly, starting 2025-01-05) written as literal
#    text - force Text number format first so Excel doesn't auto-convert the
#    strings into date serials.
$weekStartDates = @(
  "2025-01-05", "2025-01-12", "2025-01-19", "2025-01-26",
  "2025-02-02", "2025-02-09", "2025-02-16", "2025-02-23",
  "2025-03-02", "2025-03-09", "2025-03-16", "2025-03-23",
  "2025-03-30", "2025-04-06", "2025-04-13", "2025-04-20"
)

for ($i = 0; $i -lt $weekStartDates.Count; $i++) {
  $row = $i + 2
  $cell = $ws.Cells.Item($row, 2)
  $cell.NumberFormat = "@"
  $cell.Value2 = $weekStartDates[$i]
}

# 4. Week labels (column A) drop the leading zero: W01 -> W1, ... W09 -> W9
#    (W10..W16 are unchanged).
for ($row = 2; $row -le 17; $row++) {
  $week = $row - 1
  $ws.Cells.Item($row, 1).Value2 = "W$week"
}

# 5. Corrected MyForecast value for week 1 (column D after the insert).
$ws.Cells.Item(2, 4).Value2 = 102

# 6. is_holiday_week (now column J) should be a real boolean, not a number.
for ($row = 2; $row -le 17; $row++) {
  $ws.Cells.Item($row, 10).Value2 = $false
}

# 7. Update the dependent summary totals (kept as text, matching the rest
#    of the column).
$summary = $wb.Worksheets.Item("Summary")

$summaryB9 = $summary.Cells.Item(9, 2)
$summaryB9.NumberFormat = "@"
$summaryB9.Value2 = "1790"

$summaryB11 = $summary.Cells.Item(11, 2)
$summaryB11.NumberFormat = "@"
$summaryB11.Value2 = "427"
